# Auto-generated edit script applying the diff to Sheets/Diabolos_Profits.xlsx
# Each sheet corresponds to one job class tab (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ----- ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 1890
$ws.Range("I20").Value = 1890
$ws.Range("K20").Value = 1890
$ws.Range("M20").Value = -1660
$ws.Range("H35").Value = 1890
$ws.Range("I35").Value = 1890
$ws.Range("K35").Value = 1890
$ws.Range("M35").Value = -1511
$ws.Range("H53").Value = 4011.4375
$ws.Range("I53").Value = 227.11111
$ws.Range("J53").Value = 8877
$ws.Range("K53").Value = 227.11111
$ws.Range("L53").Value = 8877
$ws.Range("M53").Value = 409.88889
$ws.Range("N53").Value = -10151
$ws.Range("H62").Value = 907149.3
$ws.Range("I62").Value = 1475143
$ws.Range("J62").Value = 111958.2
$ws.Range("K62").Value = 1475143
$ws.Range("L62").Value = 111958.2
$ws.Range("M62").Value = -1474519
$ws.Range("N62").Value = -113206.2
$ws.Range("H65").Value = 907149.3
$ws.Range("I65").Value = 1475143
$ws.Range("J65").Value = 111958.2
$ws.Range("K65").Value = 7375715
$ws.Range("L65").Value = 559791
$ws.Range("M65").Value = -7372595
$ws.Range("N65").Value = -566031
$ws.Range("H86").Value = 6908553
$ws.Range("I86").Value = 7773.143
$ws.Range("K86").Value = 7773.143
$ws.Range("M86").Value = -6650.143
$ws.Range("H89").Value = 6908553
$ws.Range("I89").Value = 7773.143
$ws.Range("K89").Value = 38865.715
$ws.Range("M89").Value = -33249.715
$ws.Range("H137").Value = 1349.75
$ws.Range("I137").Value = 1349.75
$ws.Range("K137").Value = 4049.25
$ws.Range("M137").Value = -1499.25

# ----- ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 962.85486
$ws.Range("I32").Value = 962.85486
$ws.Range("K32").Value = 962.85486
$ws.Range("M32").Value = -675.85486
$ws.Range("H132").Value = 55557540
$ws.Range("I132").Value = 58825390
$ws.Range("K132").Value = 176476170
$ws.Range("M132").Value = -176473640
$ws.Range("L134").ClearContents()
$ws.Range("H134").Value = 15390
$ws.Range("I134").Value = 15390
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 15390
$ws.Range("N134").Value = 0
$ws.Range("M134").Value = -10320

# ----- BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3277.3333
$ws.Range("I86").Value = 4010.4443
$ws.Range("K86").Value = 4010.4443
$ws.Range("M86").Value = -2887.4443
$ws.Range("H89").Value = 3277.3333
$ws.Range("I89").Value = 4010.4443
$ws.Range("K89").Value = 20052.2215
$ws.Range("M89").Value = -14436.2215
$ws.Range("H99").Value = 1129.9
$ws.Range("I99").Value = 1084.2858
$ws.Range("K99").Value = 1084.2858
$ws.Range("M99").Value = 413.7141999999999

# ----- CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3061.0857
$ws.Range("I31").Value = 1385.4286
$ws.Range("J31").Value = 3480
$ws.Range("K31").Value = 1385.4286
$ws.Range("L31").Value = 3480
$ws.Range("M31").Value = -1090.4286
$ws.Range("N31").Value = -4070
$ws.Range("H34").Value = 3061.0857
$ws.Range("I34").Value = 1385.4286
$ws.Range("J34").Value = 3480
$ws.Range("K34").Value = 1385.4286
$ws.Range("L34").Value = 3480
$ws.Range("M34").Value = -1183.4286
$ws.Range("N34").Value = -3884
$ws.Range("H86").Value = 107430.375
$ws.Range("I86").Value = 210112.25
$ws.Range("J86").Value = 4748.5
$ws.Range("K86").Value = 210112.25
$ws.Range("L86").Value = 4748.5
$ws.Range("M86").Value = -208989.25
$ws.Range("N86").Value = -6994.5
$ws.Range("H89").Value = 107430.375
$ws.Range("I89").Value = 210112.25
$ws.Range("J89").Value = 4748.5
$ws.Range("K89").Value = 1050561.25
$ws.Range("L89").Value = 23742.5
$ws.Range("M89").Value = -1044945.25
$ws.Range("N89").Value = -34974.5
$ws.Range("H107").Value = 2161.5417
$ws.Range("I107").Value = 2130.5881
$ws.Range("J107").Value = 2236.7144
$ws.Range("K107").Value = 2130.5881
$ws.Range("L107").Value = 2236.7144
$ws.Range("M107").Value = -210.5880999999999
$ws.Range("N107").Value = -6076.7144

# ----- CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H47").Value = 800
$ws.Range("I47").Value = 800
$ws.Range("K47").Value = 2400
$ws.Range("M47").Value = -1969
$ws.Range("H63").Value = 15763.647
$ws.Range("I63").Value = 15763.647
$ws.Range("K63").Value = 47290.94100000001
$ws.Range("M63").Value = -46541.94100000001
$ws.Range("H66").Value = 15763.647
$ws.Range("I66").Value = 15763.647
$ws.Range("K66").Value = 141872.823
$ws.Range("M66").Value = -138128.823
$ws.Range("H101").Value = 30000
$ws.Range("J101").Value = 30000
$ws.Range("L101").Value = 90000
$ws.Range("N101").Value = -94868

# ----- GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 10018
$ws.Range("J40").Value = 10018
$ws.Range("L40").Value = 10018
$ws.Range("N40").Value = -10320
$ws.Range("H43").Value = 2724.5625
$ws.Range("I43").Value = 2724.5625
$ws.Range("K43").Value = 2724.5625
$ws.Range("M43").Value = -2573.5625

# ----- LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("L36").ClearContents()
$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("N36").Value = 0
$ws.Range("H40").Value = 2430.9167
$ws.Range("I40").Value = 2417.1
$ws.Range("J40").Value = 2500
$ws.Range("K40").Value = 2417.1
$ws.Range("L40").Value = 2500
$ws.Range("M40").Value = -2281.1
$ws.Range("N40").Value = -2772
$ws.Range("H68").Value = 19033.666
$ws.Range("I68").Value = 7299.75
$ws.Range("J68").Value = 42501.5
$ws.Range("K68").Value = 7299.75
$ws.Range("L68").Value = 42501.5
$ws.Range("M68").Value = -6550.75
$ws.Range("N68").Value = -43999.5
$ws.Range("H71").Value = 19033.666
$ws.Range("I71").Value = 7299.75
$ws.Range("J71").Value = 42501.5
$ws.Range("K71").Value = 36498.75
$ws.Range("L71").Value = 212507.5
$ws.Range("M71").Value = -32754.75
$ws.Range("N71").Value = -219995.5
$ws.Range("H93").Value = 1388.1818
$ws.Range("J93").Value = 1666.6666
$ws.Range("L93").Value = 1666.6666
$ws.Range("N93").Value = -4162.6666
$ws.Range("H132").Value = 3019.3667
$ws.Range("J132").Value = 3610.4167
$ws.Range("L132").Value = 10831.2501
$ws.Range("N132").Value = -15891.2501

# ----- WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("M51").ClearContents()
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("H81").Value = 9528296
$ws.Range("I81").Value = 1446.7142
$ws.Range("K81").Value = 2893.4284
$ws.Range("M81").Value = -1832.4284
$ws.Range("H84").Value = 9528296
$ws.Range("I84").Value = 1446.7142
$ws.Range("K84").Value = 14467.142
$ws.Range("M84").Value = -9163.142
$ws.Range("H100").Value = 5223
$ws.Range("I100").Value = 7764.7144
$ws.Range("K100").Value = 15529.4288
$ws.Range("M100").Value = -14988.4288
$ws.Range("H132").Value = 4519.5415
$ws.Range("I132").Value = 4614.467
$ws.Range("J132").Value = 4361.3335
$ws.Range("K132").Value = 13843.401
$ws.Range("L132").Value = 13084.0005
$ws.Range("M132").Value = -11313.401
$ws.Range("N132").Value = -18144.0005
